# Update (Removed Auto Arima)
# Refresh the Prophet/Amazon forecast figures on "Forecast Comparison" and
# the roll-up totals on "Summary" now that Auto-ARIMA has been dropped from
# the forecasting pipeline.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# -- "Forecast Comparison" sheet --------------------------------------------
# Columns: C = Prophet Forecast, D = Amazon Mean, E = Amazon P70,
#          F = Amazon P80, G = Amazon P90

$wsForecast.Range("D2").Value = 142
$wsForecast.Range("E2").Value = 174
$wsForecast.Range("F2").Value = 214
$wsForecast.Range("G2").Value = 280

$wsForecast.Range("D3").Value = 126
$wsForecast.Range("E3").Value = 154
$wsForecast.Range("F3").Value = 195
$wsForecast.Range("G3").Value = 264

$wsForecast.Range("D4").Value = 109
$wsForecast.Range("E4").Value = 133
$wsForecast.Range("F4").Value = 167
$wsForecast.Range("G4").Value = 220

$wsForecast.Range("D5").Value = 83
$wsForecast.Range("E5").Value = 102
$wsForecast.Range("F5").Value = 126
$wsForecast.Range("G5").Value = 166

$wsForecast.Range("D6").Value = 85
$wsForecast.Range("E6").Value = 105
$wsForecast.Range("F6").Value = 131
$wsForecast.Range("G6").Value = 173

$wsForecast.Range("D7").Value = 80
$wsForecast.Range("E7").Value = 98
$wsForecast.Range("F7").Value = 123
$wsForecast.Range("G7").Value = 163

$wsForecast.Range("C8").Value = 92
$wsForecast.Range("D8").Value = 81
$wsForecast.Range("E8").Value = 99
$wsForecast.Range("F8").Value = 125
$wsForecast.Range("G8").Value = 168

$wsForecast.Range("D9").Value = 78
$wsForecast.Range("E9").Value = 96
$wsForecast.Range("F9").Value = 121
$wsForecast.Range("G9").Value = 163

$wsForecast.Range("C10").Value = 119
$wsForecast.Range("D10").Value = 77
$wsForecast.Range("E10").Value = 94
$wsForecast.Range("F10").Value = 117
$wsForecast.Range("G10").Value = 156

$wsForecast.Range("C11").Value = 137
$wsForecast.Range("D11").Value = 75
$wsForecast.Range("E11").Value = 91
$wsForecast.Range("F11").Value = 116
$wsForecast.Range("G11").Value = 156

$wsForecast.Range("C12").Value = 134
$wsForecast.Range("D12").Value = 71
$wsForecast.Range("E12").Value = 87
$wsForecast.Range("F12").Value = 113
$wsForecast.Range("G12").Value = 157

$wsForecast.Range("D13").Value = 74
$wsForecast.Range("E13").Value = 91
$wsForecast.Range("F13").Value = 117
$wsForecast.Range("G13").Value = 160

$wsForecast.Range("D14").Value = 75
$wsForecast.Range("E14").Value = 92
$wsForecast.Range("F14").Value = 118
$wsForecast.Range("G14").Value = 161

$wsForecast.Range("D15").Value = 66
$wsForecast.Range("E15").Value = 81
$wsForecast.Range("F15").Value = 106
$wsForecast.Range("G15").Value = 148

$wsForecast.Range("D16").Value = 68
$wsForecast.Range("E16").Value = 83
$wsForecast.Range("F16").Value = 108
$wsForecast.Range("G16").Value = 149

$wsForecast.Range("D17").Value = 63
$wsForecast.Range("E17").Value = 77
$wsForecast.Range("F17").Value = 101
$wsForecast.Range("G17").Value = 140

# -- "Summary" sheet ----------------------------------------------------
# These totals are stored as text, so force a text format before writing
# the new numeric-looking strings to keep them from being re-typed as
# numbers.
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "1489"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "701"

$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "137"
